# Append the 2026-02-08 profit snapshot as a new row (row 76) to the
# bottom of the data table on the active sheet, matching the format of
# the existing rows (date stored as plain text, the rest as numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 76

# Column A holds the date as literal text (e.g. "02/07/2026" on the row
# above), not a real Excel date. Temporarily force a text number format
# so assigning the Value doesn't get auto-converted into a date serial
# number, then clear the format again so the cell ends up unstyled -
# exactly like the rest of the date column.
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "02/08/2026"
$dateCell.ClearFormats()

$ws.Cells.Item($newRow, 2).Value  = 9764.35
$ws.Cells.Item($newRow, 3).Value  = 0.2389734473272497
$ws.Cells.Item($newRow, 4).Value  = 0.7610265526727503
$ws.Cells.Item($newRow, 5).Value  = -292.51
$ws.Cells.Item($newRow, 6).Value  = -35.03
$ws.Cells.Item($newRow, 7).Value  = -23460.83
$ws.Cells.Item($newRow, 8).Value  = -75.95
$ws.Cells.Item($newRow, 9).Value  = -1019.43
$ws.Cells.Item($newRow, 10).Value = -30.4
$ws.Cells.Item($newRow, 11).Value = -24480.26
$ws.Cells.Item($newRow, 12).Value = -71.48999999999999

Write-Output "Added row $newRow (2026-02-08 snapshot)"
